$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.282.20"
$ws.Range("E2").Value = "  +5.44%  "
$ws.Range("D3").Value = "1.916.86"
$ws.Range("E3").Value = "  +5.93%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "254.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9994"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5170"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "46.09"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2979"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06838"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.01%  "
$ws.Range("D11").Value = "1.917.14"
$ws.Range("E11").Value = "  +5.95%  "
$ws.Range("E12").Value = "  +4.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07335"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6901"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "87.76"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.920"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.32%  "
$ws.Range("D17").Value = "30.288.08"
$ws.Range("E17").Value = "  +5.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008007"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +8.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9994"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("E20").Value = "  +6.62%  "
$ws.Range("D21").Value = "2.165.32"
$ws.Range("E21").Value = "  +6.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9988"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("E23").Value = "  +5.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.770"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.203"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "146.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "139.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +24.68%  "
$ws.Range("E28").Value = "  +7.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.016"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.376"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.285"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08857"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.040"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05136"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.00%  "
$ws.Range("E35").Value = "  +6.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7218"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.689"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.335"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.842"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9766"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.32%  "
$ws.Range("E41").Value = "  +6.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.225"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4337"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "105.85"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9989"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.712"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1279"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05733"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.38%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.557"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.29%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.43"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.31%  "
$ws.Range("E51").Value = "  +6.74%  "
